# Auto-generated Excel COM-interop edit script
# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: mark cells whose new text values look like plain numbers (e.g. "4.591")
# as Text format BEFORE assigning the value, so Excel does not silently convert
# them into numeric cells. This must happen before the .Value assignment below.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Step 2: write every updated cell value (coin names, links, prices, % changes).
$ws.Range("D2").Value = "19.854.50"
$ws.Range("E2").Value = "  -8.24%  "
$ws.Range("D3").Value = "1.387.39"
$ws.Range("E3").Value = "  -9.46%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "1.005"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "267.66"
$ws.Range("E6").Value = "  -7.11%  "
$ws.Range("D7").Value = "0.3611"
$ws.Range("E7").Value = "  -8.75%  "
$ws.Range("D8").Value = "0.3002"
$ws.Range("E8").Value = "  -4.52%  "
$ws.Range("D9").Value = "38.95"
$ws.Range("E9").Value = "  -7.87%  "
$ws.Range("D10").Value = "0.06320"
$ws.Range("E10").Value = "  -11.44%  "
$ws.Range("D11").Value = "0.9509"
$ws.Range("E11").Value = "  -8.61%  "
$ws.Range("D12").Value = "1.005"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "5.187"
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.395.39"
$ws.Range("E14").Value = "  -9.11%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "5.962"
$ws.Range("E15").Value = "  -9.24%  "
$ws.Range("D16").Value = "16.16"
$ws.Range("E16").Value = "  -12.26%  "
$ws.Range("D17").Value = "0.000009762"
$ws.Range("E17").Value = "  -9.82%  "
$ws.Range("D18").Value = "0.05597"
$ws.Range("E18").Value = "  -15.03%  "
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "69.07"
$ws.Range("E20").Value = "  -16.85%  "
$ws.Range("D21").Value = "5.415"
$ws.Range("E21").Value = "  -11.06%  "
$ws.Range("D22").Value = "13.99"
$ws.Range("E22").Value = "  -8.89%  "
$ws.Range("D23").Value = "10.40"
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("D24").Value = "2.246"
$ws.Range("E24").Value = "  -5.19%  "
$ws.Range("D25").Value = "19.874.95"
$ws.Range("E25").Value = "  -8.16%  "
$ws.Range("D26").Value = "2.108"
$ws.Range("E26").Value = "  -9.46%  "
$ws.Range("D27").Value = "134.20"
$ws.Range("E27").Value = "  -9.09%  "
$ws.Range("D28").Value = "16.31"
$ws.Range("E28").Value = "  -10.69%  "
$ws.Range("D29").Value = "1.552.51"
$ws.Range("E29").Value = "  -9.14%  "
$ws.Range("D30").Value = "105.97"
$ws.Range("E30").Value = "  -9.07%  "
$ws.Range("D31").Value = "3.829"
$ws.Range("E31").Value = "  -20.90%  "
$ws.Range("D32").Value = "5.135"
$ws.Range("E32").Value = "  -12.03%  "
$ws.Range("D33").Value = "0.7740"
$ws.Range("E33").Value = "  -17.28%  "
$ws.Range("D34").Value = "0.07504"
$ws.Range("E34").Value = "  -7.66%  "
$ws.Range("D35").Value = "8.161"
$ws.Range("E35").Value = "  -3.63%  "
$ws.Range("D36").Value = "1.003"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "4.591"
$ws.Range("E37").Value = "  -9.63%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.05488"
$ws.Range("E38").Value = "  -8.38%  "
$ws.Range("D39").Value = "0.1858"
$ws.Range("E39").Value = "  -7.43%  "
$ws.Range("D40").Value = "0.01971"
$ws.Range("E40").Value = "  -10.01%  "
$ws.Range("D41").Value = "1.297"
$ws.Range("E41").Value = "  -9.72%  "
$ws.Range("D42").Value = "9.814"
$ws.Range("E42").Value = "  -9.18%  "
$ws.Range("D43").Value = "1.025"
$ws.Range("E43").Value = "  -12.28%  "
$ws.Range("D44").Value = "3.446"
$ws.Range("E44").Value = "  -7.30%  "
$ws.Range("D45").Value = "0.5094"
$ws.Range("E45").Value = "  -10.87%  "
$ws.Range("D46").Value = "11.59"
$ws.Range("E46").Value = "  -10.67%  "
$ws.Range("D47").Value = "0.4875"
$ws.Range("E47").Value = "  -10.61%  "
$ws.Range("D48").Value = "106.69"
$ws.Range("E48").Value = "  -7.75%  "
$ws.Range("D49").Value = "1.686"
$ws.Range("E49").Value = "  -9.01%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "1.012"
$ws.Range("E51").Value = "  -12.71%  "

# Step 3: restore the cell style for the text-formatted cells so the saved file
# keeps the original (unstyled) look of the worksheet.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
